# Automatische test-sync: 2025-06-19 22:12:50
# Adds the new incoming mail-log entry (row 39) to the "Logs" sheet and
# bumps the matching category count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A39").Value = "Wat zijn jullie openingstijden?"
$logs.Range("B39").Value = "mailmind.test@zohomail.eu"
$logs.Range("C39").Value = "Hallo, ik zou graag willen weten wat jullie openingstijden zijn. Dank je wel!"
$logs.Range("D39").Value = "Openingstijden / Locatie"
$logs.Range("E39").Value = "Beste klant,`nDank voor uw interesse. Onze openingstijden zijn maandag t/m vrijdag van 9:00 tot 18:00 uur en zaterdag van 10:00 tot 17:00 uur. Op zondag zijn wij gesloten. Voor verdere vragen staan wij graag tot uw dienst.`nMet vriendelijke groet,`n[Naam bedrijf]"
$logs.Range("F39").Value = "2025-06-19 22:12:43"
$logs.Range("G39").Value = "Ja"

# Restore default row height (the multi-line "Antwoord" cell otherwise
# triggers an autosized row, which the source workbook does not have).
$logs.Rows.Item(39).AutoFit()

# Extend the existing conditional-formatting ranges so the new row is
# covered too (without touching rule order / dxf / priority).
$cfCategorie = $logs.Range("D2:D38").FormatConditions
for ($i = 1; $i -le $cfCategorie.Count; $i++) {
    $cfCategorie.Item($i).ModifyAppliesToRange($logs.Range("D2:D39"))
}

$cfBeantwoord = $logs.Range("G2:G38").FormatConditions
for ($i = 1; $i -le $cfBeantwoord.Count; $i++) {
    $cfBeantwoord.Item($i).ModifyAppliesToRange($logs.Range("G2:G39"))
}

$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B9").Value = 2
